$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.655.99"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "3.370.51"
$ws.Range("E3").Value = "  -4.31%  "
$ws.Range("E4").Value = "  +0.00%  "
$c1 = $ws.Range("D5")
$st1 = $c1.Style
$c1.Value = "'556.57"
$c1.Style = $st1
$ws.Range("E5").Value = "  -4.88%  "
$c2 = $ws.Range("D6")
$st2 = $c2.Style
$c2.Value = "'176.46"
$c2.Style = $st2
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").Value = "3.362.70"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("E9").Value = "  -0.05%  "
$c3 = $ws.Range("D10")
$st3 = $c3.Style
$c3.Value = "'0.629"
$c3.Style = $st3
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  -1.21%  "
$c4 = $ws.Range("D12")
$st4 = $c4.Style
$c4.Value = "'55.06"
$c4.Style = $st4
$ws.Range("E12").Value = "  -1.57%  "
$c5 = $ws.Range("D13")
$st5 = $c5.Style
$c5.Value = "'0.0000273"
$c5.Style = $st5
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("E14").Value = "  -2.60%  "
$ws.Range("D15").Value = "3.907.89"
$ws.Range("E15").Value = "  -4.32%  "
$c6 = $ws.Range("D16")
$st6 = $c6.Style
$c6.Value = "'18.42"
$c6.Style = $st6
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "3.370.05"
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("E18").Value = "  -2.78%  "
$c7 = $ws.Range("D19")
$st7 = $c7.Style
$c7.Value = "'11.86"
$c7.Style = $st7
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "64.569.43"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("E21").Value = "  -3.09%  "
$c8 = $ws.Range("D22")
$st8 = $c8.Style
$c8.Value = "'432.73"
$c8.Style = $st8
$ws.Range("E22").Value = "  +3.99%  "
$c9 = $ws.Range("D23")
$st9 = $c9.Style
$c9.Value = "'4.91"
$c9.Style = $st9
$ws.Range("E23").Value = "  +10.11%  "
$ws.Range("E24").Value = "  -5.22%  "
$c10 = $ws.Range("D25")
$st10 = $c10.Style
$c10.Value = "'84.32"
$c10.Style = $st10
$ws.Range("E25").Value = "  -1.45%  "
$c11 = $ws.Range("D26")
$st11 = $c11.Style
$c11.Value = "'13.23"
$c11.Style = $st11
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("E29").Value = "  -4.41%  "
$ws.Range("E30").Value = "  -2.35%  "
$c12 = $ws.Range("D31")
$st12 = $c12.Style
$c12.Value = "'6.65"
$c12.Style = $st12
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -2.64%  "
$c13 = $ws.Range("D33")
$st13 = $c13.Style
$c13.Value = "'577.93"
$c13.Style = $st13
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("E34").Value = "  -3.01%  "
$c14 = $ws.Range("D35")
$st14 = $c14.Style
$c14.Value = "'58.49"
$c14.Style = $st14
$ws.Range("E35").Value = "  -2.48%  "
$ws.Range("E36").Value = "  +0.16%  "
$c15 = $ws.Range("D37")
$st15 = $c15.Style
$c15.Value = "'0.143"
$c15.Style = $st15
$ws.Range("E37").Value = "  -7.93%  "
$c16 = $ws.Range("D38")
$st16 = $c16.Style
$c16.Value = "'3.51"
$c16.Style = $st16
$ws.Range("E38").Value = "  -4.25%  "
$c17 = $ws.Range("D39")
$st17 = $c17.Style
$c17.Value = "'35.82"
$c17.Style = $st17
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("E41").Value = "  -4.56%  "
$ws.Range("D42").Value = "3.118.50"
$ws.Range("E42").Value = "  -4.01%  "
$c18 = $ws.Range("D43")
$st18 = $c18.Style
$c18.Value = "'0.999"
$c18.Style = $st18
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  -5.63%  "
$c19 = $ws.Range("D45")
$st19 = $c19.Style
$c19.Value = "'3.28"
$c19.Style = $st19
$ws.Range("E45").Value = "  -2.73%  "
$c20 = $ws.Range("D46")
$st20 = $c20.Style
$c20.Value = "'0.0411"
$c20.Style = $st20
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c21 = $ws.Range("D50")
$st21 = $c21.Style
$c21.Value = "'8.31"
$c21.Style = $st21
$ws.Range("E50").Value = "  -4.43%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c22 = $ws.Range("D51")
$st22 = $c22.Style
$c22.Value = "'135.16"
$c22.Style = $st22
$ws.Range("E51").Value = "  -2.33%  "
